$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Water Quality Data")

# Update the sub-header text for the Pond/Tank column (E3) to clarify naming convention
$ws.Range("E3").Value = "Tank Name. E.g. LP1"

# Give column E an explicit width now that the text is longer
$ws.Columns.Item(5).ColumnWidth = 11.88

# Move the active selection, matching the saved workbook state
$ws.Range("E5").Select()
